$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column H (total_arqueo_ciego), shifting the existing
# total_ventas (old H) -> I and total_operaciones (old I) -> J,
# carrying their number formats (styles) along with them.
$ws.Columns.Item(8).Insert()

# Header row
$ws.Cells.Item(1, 8).Value = "total_arqueo_ciego"
$ws.Cells.Item(1, 9).Value = "total_ventas"
$ws.Cells.Item(1, 10).Value = "total_operaciones"

# Data rows: Tienda, serie, Nombre_TPV, fecha, cierre_tpv_id, cierre_tpv_desc,
# Nombre_MdP, total_arqueo_ciego, total_ventas, total_operaciones
$data = @(
    @{Row=2;  B="V2"; C="BAR";             D=45689; E=8829; G="EUROS";        H=551.3;    I=72.75;             J=9},
    @{Row=3;  B="V2"; C="BAR";             D=45689; E=8829; G="TARJETA VISA"; H=891.6;    I=202.1;             J=24},
    @{Row=4;  B="V2"; C="BAR";             D=45689; E=8830; G="EUROS";        H=422.05;   I=83.40000000000001; J=9},
    @{Row=5;  B="V2"; C="BAR";             D=45689; E=8830; G="TARJETA VISA"; H=867;      I=139.8;             J=14},
    @{Row=6;  B="V1"; C="SERVIDOR TIENDA"; D=45689; E=8828; G="EUROS";        H=1204.3;   I=806.3200000000001; J=77},
    @{Row=7;  B="V1"; C="SERVIDOR TIENDA"; D=45689; E=8828; G="TARJETA VISA"; H=2231.23;  I=2231.03;           J=154},
    @{Row=8;  B="V1"; C="SERVIDOR TIENDA"; D=45689; E=8831; G="EUROS";        H=1130.48;  I=702.76;            J=79},
    @{Row=9;  B="V1"; C="SERVIDOR TIENDA"; D=45689; E=8831; G="SMS";          H=0;        I=4.7;               J=1},
    @{Row=10; B="V1"; C="SERVIDOR TIENDA"; D=45689; E=8831; G="TARJETA VISA"; H=1801.6;   I=1809.8;            J=154},
    @{Row=11; B="V1"; C="SERVIDOR TIENDA"; D=45689; E=8833; G="TARJETA VISA"; H=1373;     I=93.2;              J=13}
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
}
